$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values: week 18 (row 19) casos 388 -> 390, week 19 (row 20) casos 2 -> 537
$ws.Range("B19").Value = 390
$ws.Range("B20").Value = 537

# Add new week 20 row (row 21)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 3
